# Weekly update: a new week of "Vega Monumental Concepción - Acelga" data is
# prepended to the table (rows 32-33), shifting all subsequent rows down by 2
# and extending the table from row 149 to row 151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 32, pushing all rows from
# 32 onward (through the former 149) down to 34..151.
$ws.Rows("32:33").Insert()

# New row 32: "Primera" quality entry for the new week (14-09-2021 / 44453)
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44453
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112009
$ws.Range("G32").Value = "Acelga"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 700
$ws.Range("M32").Value = 650
$ws.Range("N32").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O32").Value = "Región de Ñuble"
$ws.Range("P32").Value = 650
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"

# New row 33: "Segunda" quality entry for the new week
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44453
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112009
$ws.Range("G33").Value = "Acelga"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 500
$ws.Range("N33").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O33").Value = "Región de Ñuble"
$ws.Range("P33").Value = 500
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"
